$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (2024-05-03 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update price values in column D for rows 30-33
$ws.Range("D30").Value = 1956.522
$ws.Range("D31").Value = 2316.776
$ws.Range("D32").Value = 2616.068
$ws.Range("D33").Value = 3297.817
